$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 -- this shifts the existing rows 6:75 down to 7:76,
# matching the diff's net effect (dimension grows from A1:T75 to A1:T76).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44545
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 2300
$ws.Range("O6").Value = 2300
$ws.Range("P6").Value = 2300
$ws.Range("Q6").Value = "`$/kilo"
$ws.Range("R6").Value = "Región del Maule"
$ws.Range("S6").Value = 2300
$ws.Range("T6").Value = 1
